# OPCodesRev02 - ALU / Memory / Control Signal Generator pass
#
# - Makes "(A)Instruction Generator Format(A)" the active sheet/tab
#   (it was "(B)...Format(B)" before).
# - Reworks the Format(A) worked example so Source1/Source2/Dest all
#   carry the value 4 (R[4]), flowing through the existing HEX/BIN
#   formulas and the concatenated instruction in B11.
# - Swaps the "Final Decimal Instruction" cell (B12) from BIN2HEX to
#   BIN2DEC so it matches the B11->decimal pattern used elsewhere.
# - Appends a new "Past Instructions" row (26) recording the new
#   R[4]<- R[4]+R[4] example and its decimal instruction encoding.

$wb  = $excel.ActiveWorkbook
$wsA = $wb.Worksheets.Item("(A)Instruction Generator Format(A)")

# --- Source1 / Source2 / Destination1 all become 4 -------------------
$wsA.Cells.Item(5, 2).Value = 4   # B5
$wsA.Cells.Item(5, 3).Value = 4   # C5
$wsA.Cells.Item(5, 4).Value = 4   # D5

# --- Final Decimal Instruction: BIN2HEX(B11) -> BIN2DEC(B11) ---------
$wsA.Cells.Item(12, 2).Formula = "=BIN2DEC(B11)"   # B12

# --- New worked example appended as row 26 ----------------------------
$wsA.Cells.Item(26, 1).Value = "R[4]<- R[4]+R[4]"
$wsA.Cells.Item(26, 2).Value = 554172480

# --- Make Format(A) the active tab/sheet (was Format(B)) -------------
# Format(B)'s own selection (B24) is untouched/unchanged by this edit.
$wsA.Activate() | Out-Null
$wsA.Range("B24").Select() | Out-Null
